$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 898.7646999999999
$ws.Range("I19").Value = 647.0769
$ws.Range("K19").Value = 647.0769
$ws.Range("M19").Value = -472.0769
# Row 80
$ws.Range("H80").Value = 1211.1538
$ws.Range("I80").Value = 1158.8334
$ws.Range("J80").Value = 1256
$ws.Range("K80").Value = 3476.5002
$ws.Range("L80").Value = 3768
$ws.Range("M80").Value = -2478.5002
$ws.Range("N80").Value = -5764
# Row 83
$ws.Range("H83").Value = 1211.1538
$ws.Range("I83").Value = 1158.8334
$ws.Range("J83").Value = 1256
$ws.Range("K83").Value = 10429.5006
$ws.Range("L83").Value = 11304
$ws.Range("M83").Value = -5437.500599999999
$ws.Range("N83").Value = -21288
# Row 101
$ws.Range("H101").Value = 1469.7142
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 138
$ws.Range("H138").Value = 2844.182
$ws.Range("I138").Value = 2844.182
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 8532.545999999998
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -3392.545999999998
$ws.Range("N138").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4207.511
$ws.Range("I32").Value = 4025.121
$ws.Range("J32").Value = 4709.0835
$ws.Range("K32").Value = 4025.121
$ws.Range("L32").Value = 4709.0835
$ws.Range("M32").Value = -3738.121
$ws.Range("N32").Value = -5283.0835
# Row 45
$ws.Range("H45").Value = 1098
$ws.Range("I45").Value = 1096.6
$ws.Range("J45").Value = 1105
$ws.Range("K45").Value = 1096.6
$ws.Range("L45").Value = 1105
$ws.Range("M45").Value = -719.5999999999999
$ws.Range("N45").Value = -1859
# Row 122
$ws.Range("H122").Value = 3615.2964
$ws.Range("I122").Value = 3846.3635
$ws.Range("K122").Value = 11539.0905
$ws.Range("M122").Value = -9089.0905
# Row 132
$ws.Range("H132").Value = 4126.4067
$ws.Range("I132").Value = 2107.3022
$ws.Range("J132").Value = 9552.75
$ws.Range("K132").Value = 6321.9066
$ws.Range("L132").Value = 28658.25
$ws.Range("M132").Value = -3791.9066
$ws.Range("N132").Value = -33718.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 25
$ws.Range("H25").Value = 7153.8335
$ws.Range("I25").Value = 416.66666
$ws.Range("J25").Value = 9399.556
$ws.Range("K25").Value = 416.66666
$ws.Range("L25").Value = 9399.556
$ws.Range("M25").Value = -181.66666
$ws.Range("N25").Value = -9869.556
# Row 86
$ws.Range("H86").Value = 1494.6154
$ws.Range("I86").Value = 1241.1428
$ws.Range("K86").Value = 1241.1428
$ws.Range("M86").Value = -118.1428000000001
# Row 89
$ws.Range("H89").Value = 1494.6154
$ws.Range("I89").Value = 1241.1428
$ws.Range("K89").Value = 6205.714
$ws.Range("M89").Value = -589.7139999999999
# Row 105
$ws.Range("H105").Value = 142860340
$ws.Range("I105").Value = 166669630
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 166669630
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = -166667883
$ws.Range("N105").Value = -7994

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 1100
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 1200
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 1200
$ws.Range("M23").Value = -760
$ws.Range("N23").Value = -1680
# Row 27
$ws.Range("H27").Value = 1100
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1200
$ws.Range("M27").Value = -808
$ws.Range("N27").Value = -1584
# Row 31
$ws.Range("H31").Value = 2177227.2
$ws.Range("I31").Value = 1670.28
$ws.Range("J31").Value = 4767176
$ws.Range("K31").Value = 1670.28
$ws.Range("L31").Value = 4767176
$ws.Range("M31").Value = -1375.28
$ws.Range("N31").Value = -4767766
# Row 34
$ws.Range("H34").Value = 2177227.2
$ws.Range("I34").Value = 1670.28
$ws.Range("J34").Value = 4767176
$ws.Range("K34").Value = 1670.28
$ws.Range("L34").Value = 4767176
$ws.Range("M34").Value = -1468.28
$ws.Range("N34").Value = -4767580
# Row 86
$ws.Range("H86").Value = 5255.3335
$ws.Range("I86").Value = 4351
$ws.Range("J86").Value = 5901.2856
$ws.Range("K86").Value = 4351
$ws.Range("L86").Value = 5901.2856
$ws.Range("M86").Value = -3228
$ws.Range("N86").Value = -8147.2856
# Row 89
$ws.Range("H89").Value = 5255.3335
$ws.Range("I89").Value = 4351
$ws.Range("J89").Value = 5901.2856
$ws.Range("K89").Value = 21755
$ws.Range("L89").Value = 29506.428
$ws.Range("M89").Value = -16139
$ws.Range("N89").Value = -40738.428
# Row 107
$ws.Range("H107").Value = 2809.2727
$ws.Range("I107").Value = 1938.75
$ws.Range("J107").Value = 5130.6665
$ws.Range("K107").Value = 1938.75
$ws.Range("L107").Value = 5130.6665
$ws.Range("M107").Value = -18.75
$ws.Range("N107").Value = -8970.666499999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 834.3333
$ws.Range("I14").Value = 834.3333
$ws.Range("K14").Value = 2502.9999
$ws.Range("M14").Value = -2329.9999
# Row 33
$ws.Range("H33").Value = 85.30768999999999
$ws.Range("I33").Value = 57.5
$ws.Range("K33").Value = 345
$ws.Range("M33").Value = -62
# Row 34
$ws.Range("H34").Value = 1901.9445
$ws.Range("J34").Value = 3281.3
$ws.Range("L34").Value = 9843.900000000001
$ws.Range("N34").Value = -10011.9
# Row 39
$ws.Range("H39").Value = 4802.1
$ws.Range("J39").Value = 9599.799999999999
$ws.Range("L39").Value = 28799.4
$ws.Range("N39").Value = -29387.4
# Row 48
$ws.Range("H48").Value = 1180
$ws.Range("I48").Value = 700
$ws.Range("J48").Value = 1900
$ws.Range("K48").Value = 2100
$ws.Range("L48").Value = 5700
$ws.Range("M48").Value = -1850
$ws.Range("N48").Value = -6200
# Row 55
$ws.Range("H55").Value = 2311.2856
$ws.Range("J55").Value = 7305.5
$ws.Range("L55").Value = 21916.5
$ws.Range("N55").Value = -22270.5
# Row 137
$ws.Range("H137").Value = 1874.1333
$ws.Range("I137").Value = 1190.3334
$ws.Range("J137").Value = 2899.8333
$ws.Range("K137").Value = 3571.0002
$ws.Range("L137").Value = 8699.499899999999
$ws.Range("M137").Value = 1528.9998
$ws.Range("N137").Value = -18899.4999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 48000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 48000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 104
$ws.Range("H104").Value = 80000
$ws.Range("J104").Value = 80000
$ws.Range("L104").Value = 80000
$ws.Range("N104").Value = -86988
# Row 126
$ws.Range("H126").Value = 2727.2856
$ws.Range("I126").Value = 1488.1333
$ws.Range("K126").Value = 4464.3999
$ws.Range("M126").Value = -1994.3999
# Row 132
$ws.Range("H132").Value = 6256.07
$ws.Range("I132").Value = 4787.0835
$ws.Range("K132").Value = 14361.2505
$ws.Range("M132").Value = -11831.2505

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2954.0732
$ws.Range("I22").Value = 1908.2916
$ws.Range("J22").Value = 4430.4707
$ws.Range("K22").Value = 1908.2916
$ws.Range("L22").Value = 4430.4707
$ws.Range("M22").Value = -1613.2916
$ws.Range("N22").Value = -5020.4707
# Row 27
$ws.Range("H27").Value = 2954.0732
$ws.Range("I27").Value = 1908.2916
$ws.Range("J27").Value = 4430.4707
$ws.Range("K27").Value = 1908.2916
$ws.Range("L27").Value = 4430.4707
$ws.Range("M27").Value = -1801.2916
$ws.Range("N27").Value = -4644.4707
# Row 55
$ws.Range("H55").Value = 1405.5454
$ws.Range("I55").Value = 1062.4445
$ws.Range("K55").Value = 1062.4445
$ws.Range("M55").Value = -889.4445000000001
# Row 68
$ws.Range("H68").Value = 2877.889
$ws.Range("I68").Value = 2814.4285
$ws.Range("J68").Value = 3100
$ws.Range("K68").Value = 2814.4285
$ws.Range("L68").Value = 3100
$ws.Range("M68").Value = -2065.4285
$ws.Range("N68").Value = -4598
# Row 71
$ws.Range("H71").Value = 2877.889
$ws.Range("I71").Value = 2814.4285
$ws.Range("J71").Value = 3100
$ws.Range("K71").Value = 14072.1425
$ws.Range("L71").Value = 15500
$ws.Range("M71").Value = -10328.1425
$ws.Range("N71").Value = -22988
# Row 82
$ws.Range("H82").Value = 2143.75
$ws.Range("I82").Value = 1518.1818
$ws.Range("K82").Value = 1518.1818
$ws.Range("M82").Value = -1157.1818
# Row 85
$ws.Range("H85").Value = 2143.75
$ws.Range("I85").Value = 1518.1818
$ws.Range("K85").Value = 1518.1818
$ws.Range("M85").Value = -270.1818000000001
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1660.3334
$ws.Range("I96").Value = 1245.5
$ws.Range("K96").Value = 1245.5
$ws.Range("M96").Value = 127.5
# Row 117
$ws.Range("H117").Value = 80409
$ws.Range("J117").Value = 80409
$ws.Range("L117").Value = 80409
$ws.Range("N117").Value = -89587
# Row 126
$ws.Range("H126").Value = 60130.445
$ws.Range("I126").Value = 72016.53
$ws.Range("K126").Value = 216049.59
$ws.Range("M126").Value = -213579.59
# Row 132
$ws.Range("H132").Value = 5022.0566
$ws.Range("I132").Value = 3085.121
$ws.Range("K132").Value = 9255.363000000001
$ws.Range("M132").Value = -6725.363000000001
